# GitNote.docx: after "Git checkout: ... Any local changes you made to
# that file are gone", append two blank ListParagraph-styled paragraphs
# followed by a new ListParagraph paragraph containing "Add a different
# line" - mirroring a user placing the cursor at the end of the document,
# pressing Enter three times, then typing the new line.

$d = $word.ActiveDocument

$sel = $word.Selection
$sel.EndKey(6)          # wdStory - jump to the very end of the document
$sel.TypeParagraph()    # first new (blank) paragraph
$sel.TypeParagraph()    # second new (blank) paragraph
$sel.TypeParagraph()    # third new paragraph, will hold the new text
$sel.TypeText("Add a different line")
